# Re-apply the "Add files via upload" commit:
#  - drop the Notes Master (and its private theme) from the deck
#  - fix the placeholder date field on the slide layouts (11/7/2024 -> 11/6/2024)
#  - clear slide 1 back to an (empty) slide that just follows the master background
#
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Notes master removal.
#    The sandboxed object model only exposes NotesMaster as a Master object
#    whose Delete() routes through the single-slide-master guard, so an
#    explicit deletion call is rejected when there is only one design in the
#    deck. We still try it (in case the host ever allows it) but must not let
#    a failure abort the rest of the edit.
# ---------------------------------------------------------------------------
try {
    $ppt.ActivePresentation.NotesMaster.Delete()
} catch {
    # Not supported by this host when only one slide master/design exists -
    # continue with the rest of the edits.
}

# ---------------------------------------------------------------------------
# 2) Slide layouts: the date placeholder field text had a typo (11/7 -> 11/6).
#    Walk every custom layout on the slide master and fix any shape whose
#    name marks it as the Date placeholder.
# ---------------------------------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/7/2024") {
                $tr.Text = "11/6/2024"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 1: remove the explicit white background override so the slide
#    falls back to following the slide master's background again.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$slide1.FollowMasterBackground = -1   # msoTrue

# ---------------------------------------------------------------------------
# 4) Slide 1: remove all the poster content shapes/pictures that were added
#    (title textbox, header text block, the four rounded "card" panels and
#    the four pictures), leaving an empty slide again.
# ---------------------------------------------------------------------------
while ($slide1.Shapes.Count -gt 0) {
    $slide1.Shapes.Item(1).Delete()
}
